# IndicatorInformation.xlsx touch-up before going on leave:
#  - fix the "Learning environment" category label's capitalisation
#    (it becomes "Learning Environment" for the 7 rows that use it)
#  - fix the "Scottosh" typo in the school-statistics source label
#  - turn on AutoFilter for the data table (A1:H71)
#  - leave the active cell on H7 (where work was interrupted)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_indicators")

# --- Content fixes -------------------------------------------------

# Capitalisation fix: "Learning environment" -> "Learning Environment"
# applies to every row in the "Learning Environment" determinant block (rows 42-48)
foreach ($r in 42..48) {
    $ws.Cells.Item($r, 3).Value = "Learning Environment"
}

# Typo fix: "Scottosh Government School Statistics" -> "Scottish Government School Statistics"
$ws.Cells.Item(42, 8).Value = "Scottish Government School Statistics"

# --- Turn on AutoFilter for the table -------------------------------

$ws.Range("A1:H71").AutoFilter() | Out-Null

try {
    $filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=all_indicators!`$A`$1:`$H`$71")
    $filterName.Visible = $false
} catch {
    # already defined (e.g. re-running against an already-filtered sheet) - ignore
}

# --- Restore the cursor position the author left it at --------------

$ws.Range("H7").Select() | Out-Null
